# 16-1-1-1.xlsx — "Add files via upload" edit
#
# Content changes applied:
#  1. B10 ("Сайт организации"): "www.stat.kg" -> "www.stat.gov.kg", and the
#     cell becomes a live hyperlink pointing at that site.
#  2. B1 ("Индикатор" code): "16.1.1.1a" -> "16.1.1.1".
#  3. Active selection moves from G6 to B1 and the saved scroll position
#     resets to the top of the sheet.
#
# NOTE on ordering: B10 is updated *before* B1 so that the shared-string
# table gets the two new strings appended in the same order the author's
# workbook shows them ("www.stat.gov.kg" then "16.1.1.1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Organization website text + hyperlink
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Hyperlinks.Add($ws.Range("B10"), "http://www.stat.gov.kg")

# 2) Indicator code correction
$ws.Range("B1").Value = "16.1.1.1"

# 3) Reset selection/view to B1
$ws.Range("B1").Select()
